$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update the "time_taken" timestamps (column F) on the "data" sheet ---
$ws1.Range("F2").Value  = "2021-10-05 14:35:13.620533"
$ws1.Range("F3").Value  = "2021-10-05 14:35:13.620541"
$ws1.Range("F4").Value  = "2021-10-05 14:35:13.620544"
$ws1.Range("F5").Value  = "2021-10-05 14:35:13.620547"
$ws1.Range("F6").Value  = "2021-10-05 14:35:13.620550"
$ws1.Range("F7").Value  = "2021-10-05 14:35:13.620553"
$ws1.Range("F8").Value  = "2021-10-05 14:35:13.620556"
$ws1.Range("F9").Value  = "2021-10-05 14:35:13.620558"
$ws1.Range("F10").Value = "2021-10-05 14:35:13.620561"
$ws1.Range("F11").Value = "2021-10-05 14:35:13.620564"
$ws1.Range("F12").Value = "2021-10-05 14:35:13.620566"
$ws1.Range("F13").Value = "2021-10-05 14:35:13.620569"
$ws1.Range("F14").Value = "2021-10-05 14:35:13.620571"
$ws1.Range("F15").Value = "2021-10-05 14:35:13.620574"
$ws1.Range("F16").Value = "2021-10-05 14:35:13.620576"
$ws1.Range("F17").Value = "2021-10-05 14:35:13.620579"
$ws1.Range("F18").Value = "2021-10-05 14:35:13.620582"

# --- Add a new "metadata" sheet after "data" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "metadata"

# Header row (B1:G1) - reuse the bold/bordered header style from the "data" sheet
$ws1.Range("B1:F1").Copy() | Out-Null
$ws2.Range("B1:F1").PasteSpecial(-4122) | Out-Null
$ws1.Range("B1").Copy() | Out-Null
$ws2.Range("G1").PasteSpecial(-4122) | Out-Null

$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row (A2:G2) - reuse the numeric index style from the "data" sheet for A2
$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A2").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Value = 0

$ws2.Range("B2").Value = "Pharmacogenomics_Paediatric"
$ws2.Range("C2").Value = 3271

$d2 = $ws2.Cells.Item(2, 4)
$d2.NumberFormat = "@"
$d2.Value = "0.50"
$d2.Style = "Normal"

$ws2.Range("E2").Value = "2020-08-27T10:53:11.205850Z"
$ws2.Range("F2").Value = "2021-10-05 14:35:13.616799"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3271/?format=json"

Write-Output "metadata sheet added and timestamps refreshed"
